# Applies the "merge main & resolve conflicts" edit to
# StructureDefinition-FrRangeUcum.xlsx:
#   - Metadata sheet: URL, Date updated; Copyright value cleared.
#   - Elements sheet: Range.low (row 5) and Range.high (row 6) rows get
#     updated "Is Summary?", "Comments", "Condition(s)", "Constraint(s)",
#     "Mapping: HL7 v2 Mapping" and "Mapping: RIM Mapping" values.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "https://hl7.fr/fhir/fr/medication/StructureDefinition/FrRangeUcum"
$meta.Range("B8").Value = "2024-12-26T10:27:36+00:00"
$meta.Range("B14").Value = ""

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 5 = Range.low
$elements.Range("J5").Value = "Y"
$elements.Range("N5").Value = "If the low element is missing, the low boundary is not known."
$elements.Range("AI5").Value = ""
$elements.Range("AJ5").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
"
$elements.Range("AK5").Value = "NR.1"
$elements.Range("AL5").Value = "./low"

# Row 6 = Range.high
$elements.Range("J6").Value = "Y"
$elements.Range("N6").Value = "If the high element is missing, the high boundary is not known."
$elements.Range("AI6").Value = ""
$elements.Range("AJ6").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
"
$elements.Range("AK6").Value = "NR.2"
$elements.Range("AL6").Value = "./high"
